$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.020.25'
$ws.Range('E2').Value = '  +1.70%  '
$ws.Range('D3').Value = '1.863.44'
$ws.Range('E3').Value = '  +0.97%  '
$ws.Range('E4').Value = '  -0.36%  '
$ws.Range('D5').Value = "'336.02"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.69%  '
$ws.Range('E6').Value = '  -0.29%  '
$ws.Range('D7').Value = "'0.4705"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.48%  '
$ws.Range('D8').Value = "'0.3902"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.35%  '
$ws.Range('D9').Value = "'46.84"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.02%  '
$ws.Range('D10').Value = "'0.07963"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.01%  '
$ws.Range('D11').Value = "'0.9849"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.27%  '
$ws.Range('D12').Value = "'21.52"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.42%  '
$ws.Range('D13').Value = "'5.958"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.08%  '
$ws.Range('D14').Value = '1.853.75'
$ws.Range('E14').Value = '  -0.07%  '
$ws.Range('D15').Value = "'7.210"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.09%  '
$ws.Range('D16').Value = "'91.41"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.31%  '
$ws.Range('D17').Value = "'1.003"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.26%  '
$ws.Range('D18').Value = "'0.00001043"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.84%  '
$ws.Range('D19').Value = "'0.06615"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.45%  '
$ws.Range('D20').Value = "'17.56"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.46%  '
$ws.Range('E21').Value = '  -0.32%  '
$ws.Range('D22').Value = '28.029.54'
$ws.Range('E22').Value = '  +1.77%  '
$ws.Range('D23').Value = "'5.408"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.42%  '
$ws.Range('E24').Value = '  +1.08%  '
$ws.Range('D25').Value = "'2.291"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.47%  '
$ws.Range('D26').Value = '2.069.33'
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('D27').Value = "'159.26"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.42%  '
$ws.Range('D28').Value = "'19.51"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.41%  '
$ws.Range('D29').Value = "'2.108"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.24%  '
$ws.Range('D30').Value = "'5.481"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.56%  '
$ws.Range('D31').Value = "'119.38"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.33%  '
$ws.Range('D32').Value = "'0.9648"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.18%  '
$ws.Range('D33').Value = "'0.09488"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.98%  '
$ws.Range('D34').Value = "'3.577"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.31%  '
$ws.Range('D35').Value = "'5.313"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.22%  '
$ws.Range('D37').Value = "'0.02262"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.59%  '
$ws.Range('E38').Value = '  +1.03%  '
$ws.Range('D39').Value = "'8.319"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.53%  '
$ws.Range('D40').Value = "'1.166"
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Value = "'1.001"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.25%  '
$ws.Range('D42').Value = "'0.5934"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.74%  '
$ws.Range('D43').Value = "'0.1873"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.62%  '
$ws.Range('D44').Value = "'10.24"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.92%  '
$ws.Range('D45').Value = "'1.296"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.12%  '
$ws.Range('D46').Value = "'0.5581"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('D47').Value = "'12.14"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.08%  '
$ws.Range('D48').Value = "'1.963"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.87%  '
$ws.Range('D49').Value = "'0.06874"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.59%  '
$ws.Range('D50').Value = "'111.64"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.56%  '
$ws.Range('E51').Value = '  -32.71%  '
